# "added harvester and experiment design"
# Sheet1 gets three new columns of per-sample metadata filled in across rows 2-5:
#   B (harvester)        -> S.GISH
#   F (strain)            -> KN99alpha
#   D (experimentDesign)  -> 90minuteInduction
#
# Order matters for matching the shared-string table layout produced by Excel
# (new unique strings are appended in first-seen order): B, then F, then D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B5").Value = "S.GISH"
$ws.Range("F2:F5").Value = "KN99alpha"
$ws.Range("D2:D5").Value = "90minuteInduction"

# Column widths picked up by Excel's "best fit" after the new text landed.
$ws.Columns.Item(4).ColumnWidth = 14.6666666666667
$ws.Columns.Item(5).ColumnWidth = 19.6666666666667

# Matches the saved selection left on the new experimentDesign column.
$ws.Range("D3:D5").Select()
